$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest values.
# Price column values are forced to text (format "@") so that numeric-looking
# strings (e.g. "222.46") are not auto-converted to numbers by Excel, then the
# style is reset back to Normal so no stray cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.124.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.795.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.549'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.96'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.283'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0715'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.96%  '
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.052.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.794.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.115.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0781'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  +2.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.69%  '
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.414.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  -1.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.942'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.46%  '
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0494'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.951.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0121'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.30%  '
